$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 50 contains the "[No abstract available]" entry that needs to be removed.
$ws.Rows.Item(50).Delete()

# Update selection/top-left to match what Excel would leave after deleting the row.
$ws.Range("A50").Select()
$excel.ActiveWindow.ScrollRow = 49
